# Price update for 2026-02-07
# Appends one new row (Date, Price, Discount, Incredible) to the tracking
# sheet, mirroring the existing data rows' layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty row right after the existing data block (column A).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$rowValues = @("2026-02-07", "81400000", "0", "0")

# Every existing cell in the sheet is stored as literal text (shared
# strings), even though several columns look numeric/date-like. Pre-format
# the new cells as Text before assigning so Excel doesn't auto-convert the
# date/number-looking strings into a real date serial or numeric value.
$targetRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 4))
$targetRange.NumberFormat = "@"

for ($col = 1; $col -le $rowValues.Length; $col++) {
    $ws.Cells.Item($newRow, $col).Value = $rowValues[$col - 1]
}

# Restore the default "Normal" style so the new row matches the formatting
# (i.e. lack thereof) of every other row in the sheet.
$targetRange.Style = "Normal"
